$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.840.27'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.77%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.632.85'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.25%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.16'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.52%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.47'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.39%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.595'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.57%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.53'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.54%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.71%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.73%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.377'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.89%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.101.93'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.34%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.09'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.05%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.872.63'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.89%  '

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.86%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.646.66'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.58%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.55'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.26%  '

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.34%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '350.91'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.06%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.83%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.36%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.526'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.07%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.60'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.67%  '

# Row 25
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.50'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +5.41%  '

# Row 26
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.162'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.85%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.993'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.54%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.58%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0810'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.30%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.81'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.83%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.80'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.72%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.08%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.91'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.07%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +10.35%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.25%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.34%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '337.56'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +9.47%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.13'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.73%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.907'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +7.93%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.10'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.08%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.30'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.79%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '135.77'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.15%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0573'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.39%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.88'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.85%  '

# Row 46
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.623'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.26%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.23'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.50%  '

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.50%  '

# Row 49
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0991'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.12%  '

# Row 50
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.997'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.02%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.77'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.48%  '
